$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.053.51'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '2.049.11'
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("E6").Value = '  +1.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.18'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.00%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.380'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0778'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.21%  '
$ws.Range("E11").Value = '  +1.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.70'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.41%  '
$ws.Range("D13").Value = '2.347.94'
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.806'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.22%  '
$ws.Range("D16").Value = '2.045.62'
$ws.Range("E16").Value = '  -0.47%  '
$ws.Range("D17").Value = '37.080.07'
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +16.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.26%  '
$ws.Range("E20").Value = '  -0.79%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '236.03'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.37%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  -2.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.46%  '
$ws.Range("E27").Value = '  +0.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.17%  '
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.85%  '
$ws.Range("E31").Value = '  +2.51%  '
$ws.Range("E32").Value = '  -2.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.41'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.28%  '
$ws.Range("E34").Value = '  +1.83%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  -1.94%  '
$ws.Range("E37").Value = '  -1.26%  '
$ws.Range("E38").Value = '  +4.68%  '
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("E40").Value = '  +11.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.08'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +27.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0221'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.31'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.19%  '
$ws.Range("E44").Value = '  -1.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '95.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.71%  '
$ws.Range("E46").Value = '  +1.90%  '
$ws.Range("D47").Value = '1.274.75'
$ws.Range("E47").Value = '  -1.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.86'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.25%  '
$ws.Range("D49").Value = '2.236.71'
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("E50").Value = '  -1.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.40'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -22.38%  '
